$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 15: add the "gameOverviewInitialize" key with its explanatory text.
# Write the Value-column string first so it lands at shared-string index 25
# and the Key-column string second so it lands at index 26 (matches target).
$ws.Range("B15").Value = "true: die rundenübersicht wird für dieses Spiel neu aufgerufen (aus hauptmenü heraus)`nfalse: die rundenübersicht wird lediglich aktualisiert!"
$ws.Range("A15").Value = "gameOverviewInitialize"

# New cell style for B15: wrap text.
$ws.Range("B15").WrapText = $true

# Row height for the wrapped row.
$ws.Rows.Item(15).RowHeight = 45

# Grow the table to include the new row.
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.Resize($ws.Range("A3:B15"))

# Move the active selection to B12.
[void]$ws.Range("B12").Select()
